$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44530
$ws.Range("K2").Value = 'Castle Brite'
$ws.Range("L2").Value = 'Especial'
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 30000
$ws.Range("O2").Value = 30000
$ws.Range("P2").Value = 30000
$ws.Range("Q2").Value = '$/caja 18 kilos'
$ws.Range("R2").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S2").Value = 1667
$ws.Range("T2").Value = 18

# Row 3
$ws.Range("D3").Value = 44530
$ws.Range("K3").Value = 'Castle Brite'
$ws.Range("M3").Value = 100
$ws.Range("Q3").Value = '$/bandeja 10 kilos'
$ws.Range("R3").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S3").Value = 2000
$ws.Range("T3").Value = 10

# Row 4
$ws.Range("D4").Value = 44530
$ws.Range("K4").Value = 'Castle Brite'
$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 18000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 18000
$ws.Range("Q4").Value = '$/bandeja 10 kilos'
$ws.Range("R4").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S4").Value = 1800
$ws.Range("T4").Value = 10

# Row 5
$ws.Range("D5").Value = 44530
$ws.Range("K5").Value = 'Castle Brite'
$ws.Range("L5").Value = 'Tercera'
$ws.Range("M5").Value = 50
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = '$/bandeja 10 kilos'
$ws.Range("R5").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S5").Value = 1500
$ws.Range("T5").Value = 10

# Row 6
$ws.Range("D6").Value = 44931
$ws.Range("K6").Value = 'Dina'
$ws.Range("M6").Value = 55
$ws.Range("N6").Value = 17000
$ws.Range("O6").Value = 17000
$ws.Range("P6").Value = 17000
$ws.Range("Q6").Value = '$/bandeja 10 kilos'
$ws.Range("S6").Value = 1700
$ws.Range("T6").Value = 10

# Row 7
$ws.Range("D7").Value = 44931
$ws.Range("K7").Value = 'Modesto'
$ws.Range("M7").Value = 45
$ws.Range("Q7").Value = '$/caja 15 kilos'
$ws.Range("R7").Value = 'Región de O''Higgins'
$ws.Range("S7").Value = 1333
$ws.Range("T7").Value = 15

# Row 8
$ws.Range("D8").Value = 44539
$ws.Range("M8").Value = 75
$ws.Range("P8").Value = 19067
$ws.Range("S8").Value = 1059

# Row 9
$ws.Range("D9").Value = 44162
$ws.Range("M9").Value = 75
$ws.Range("P9").Value = 18933
$ws.Range("S9").Value = 1262

# Row 10
$ws.Range("D10").Value = 44162
$ws.Range("M10").Value = 55
$ws.Range("N10").Value = 23000
$ws.Range("O10").Value = 23000
$ws.Range("P10").Value = 23000
$ws.Range("Q10").Value = '$/caja 18 kilos'
$ws.Range("R10").Value = 'Provincia de Limarí'
$ws.Range("S10").Value = 1278
$ws.Range("T10").Value = 18

# Row 11
$ws.Range("D11").Value = 44554
$ws.Range("K11").Value = 'Dina'
$ws.Range("L11").Value = 'Primera'
$ws.Range("M11").Value = 35
$ws.Range("N11").Value = 20000
$ws.Range("O11").Value = 20000
$ws.Range("P11").Value = 20000
$ws.Range("Q11").Value = '$/caja 18 kilos'
$ws.Range("R11").Value = 'Provincia de Quillota'
$ws.Range("S11").Value = 1111
$ws.Range("T11").Value = 18

# Row 12
$ws.Range("D12").Value = 44935
$ws.Range("K12").Value = 'Modesto'
$ws.Range("L12").Value = 'Especial'
$ws.Range("M12").Value = 500
$ws.Range("N12").Value = 18000
$ws.Range("O12").Value = 18000
$ws.Range("P12").Value = 18000
$ws.Range("Q12").Value = '$/bandeja 10 kilos'
$ws.Range("R12").Value = 'Región de O''Higgins'
$ws.Range("S12").Value = 1800
$ws.Range("T12").Value = 10

# Row 13
$ws.Range("D13").Value = 44921
$ws.Range("K13").Value = 'Dina'
$ws.Range("M13").Value = 45
$ws.Range("Q13").Value = '$/caja 15 kilos'
$ws.Range("S13").Value = 1333
$ws.Range("T13").Value = 15

# Row 14
$ws.Range("D14").Value = 44533
$ws.Range("K14").Value = 'Castle Brite'
$ws.Range("M14").Value = 65
$ws.Range("N14").Value = 1500
$ws.Range("O14").Value = 1500
$ws.Range("P14").Value = 1500
$ws.Range("Q14").Value = '$/bandeja 10 kilos'
$ws.Range("R14").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S14").Value = 150
$ws.Range("T14").Value = 10

# Row 15
$ws.Range("D15").Value = 44902
$ws.Range("K15").Value = 'Castle Brite'
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 24000
$ws.Range("O15").Value = 24000
$ws.Range("P15").Value = 24000
$ws.Range("Q15").Value = '$/bandeja 18 kilos'
$ws.Range("R15").Value = 'Región de O''Higgins'
$ws.Range("S15").Value = 1333
$ws.Range("T15").Value = 18

# Row 16
$ws.Range("D16").Value = 44172
$ws.Range("M16").Value = 80
$ws.Range("N16").Value = 22000
$ws.Range("O16").Value = 22000
$ws.Range("P16").Value = 22000
$ws.Range("Q16").Value = '$/bandeja 18 kilos'
$ws.Range("R16").Value = 'Región de O''Higgins'
$ws.Range("S16").Value = 1222
$ws.Range("T16").Value = 18

# Row 17
$ws.Range("D17").Value = 44172
$ws.Range("M17").Value = 65
$ws.Range("N17").Value = 23000
$ws.Range("O17").Value = 23000
$ws.Range("P17").Value = 23000
$ws.Range("R17").Value = 'Región de O''Higgins'
$ws.Range("S17").Value = 1278

# Row 18
$ws.Range("D18").Value = 44552
$ws.Range("K18").Value = 'Dina'
$ws.Range("M18").Value = 55
$ws.Range("N18").Value = 20000
$ws.Range("O18").Value = 22000
$ws.Range("P18").Value = 21091
$ws.Range("Q18").Value = '$/caja 18 kilos'
$ws.Range("R18").Value = 'Provincia de Quillota'
$ws.Range("S18").Value = 1172
$ws.Range("T18").Value = 18

# Row 19
$ws.Range("D19").Value = 44551
$ws.Range("M19").Value = 45
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("Q19").Value = '$/bandeja 18 kilos'
$ws.Range("R19").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S19").Value = 1000
$ws.Range("T19").Value = 18

# Row 20
$ws.Range("D20").Value = 44159
$ws.Range("M20").Value = 85
$ws.Range("N20").Value = 20000
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 20000
$ws.Range("Q20").Value = '$/caja 15 kilos'
$ws.Range("S20").Value = 1333
$ws.Range("T20").Value = 15

# Row 22
$ws.Range("D22").Value = 44168
$ws.Range("M22").Value = 450
$ws.Range("N22").Value = 22000
$ws.Range("O22").Value = 23000
$ws.Range("P22").Value = 22444
$ws.Range("R22").Value = 'Región de O''Higgins'
$ws.Range("S22").Value = 1247

# Row 23
$ws.Range("D23").Value = 44911
$ws.Range("K23").Value = 'Castle Brite'
$ws.Range("M23").Value = 215
$ws.Range("N23").Value = 20000
$ws.Range("O23").Value = 21000
$ws.Range("P23").Value = 20419
$ws.Range("Q23").Value = '$/caja 15 kilos'
$ws.Range("R23").Value = 'Región de O''Higgins'
$ws.Range("S23").Value = 1361
$ws.Range("T23").Value = 15

# Row 24
$ws.Range("D24").Value = 44536
$ws.Range("M24").Value = 112
$ws.Range("N24").Value = 13000
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 13000
$ws.Range("Q24").Value = '$/bandeja 10 kilos granel'
$ws.Range("R24").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S24").Value = 1300

# Row 25
$ws.Range("D25").Value = 44559
$ws.Range("K25").Value = 'Modesto'
$ws.Range("M25").Value = 95
$ws.Range("N25").Value = 18000
$ws.Range("O25").Value = 18000
$ws.Range("P25").Value = 18000
$ws.Range("R25").Value = 'Provincia de Quillota'
$ws.Range("S25").Value = 1000

# Row 26
$ws.Range("D26").Value = 44923
$ws.Range("K26").Value = 'Dina'
$ws.Range("M26").Value = 45
$ws.Range("N26").Value = 20000
$ws.Range("O26").Value = 20000
$ws.Range("P26").Value = 20000
$ws.Range("Q26").Value = '$/caja 15 kilos'
$ws.Range("S26").Value = 1333
$ws.Range("T26").Value = 15

# Row 27
$ws.Range("D27").Value = 44923
$ws.Range("K27").Value = 'Dina'
$ws.Range("M27").Value = 45
$ws.Range("N27").Value = 20000
$ws.Range("O27").Value = 20000
$ws.Range("P27").Value = 20000
$ws.Range("Q27").Value = '$/caja 15 kilos'
$ws.Range("S27").Value = 1333
$ws.Range("T27").Value = 15

# Row 28
$ws.Range("D28").Value = 44937
$ws.Range("K28").Value = 'Modesto'
$ws.Range("M28").Value = 500
$ws.Range("N28").Value = 18000
$ws.Range("O28").Value = 18000
$ws.Range("P28").Value = 18000
$ws.Range("Q28").Value = '$/bandeja 10 kilos'
$ws.Range("R28").Value = 'Región de O''Higgins'
$ws.Range("S28").Value = 1800

# Row 29
$ws.Range("D29").Value = 44910
$ws.Range("K29").Value = 'Dina'
$ws.Range("L29").Value = 'Especial'
$ws.Range("M29").Value = 125
$ws.Range("N29").Value = 17000
$ws.Range("O29").Value = 17000
$ws.Range("P29").Value = 17000
$ws.Range("R29").Value = 'Región de O''Higgins'
$ws.Range("S29").Value = 1700

# Row 30
$ws.Range("D30").Value = 44910
$ws.Range("K30").Value = 'Dina'
$ws.Range("L30").Value = 'Primera'
$ws.Range("M30").Value = 110
$ws.Range("N30").Value = 22000
$ws.Range("O30").Value = 22000
$ws.Range("P30").Value = 22000
$ws.Range("S30").Value = 1222

# Row 31
$ws.Range("D31").Value = 44174
$ws.Range("L31").Value = 'Primera'
$ws.Range("M31").Value = 120
$ws.Range("N31").Value = 20000
$ws.Range("P31").Value = 21083
$ws.Range("Q31").Value = '$/bandeja 18 kilos'
$ws.Range("S31").Value = 1171

# Row 32
$ws.Range("D32").Value = 44179
$ws.Range("K32").Value = 'Modesto'
$ws.Range("M32").Value = 200
$ws.Range("N32").Value = 22000
$ws.Range("O32").Value = 22000
$ws.Range("P32").Value = 22000
$ws.Range("Q32").Value = '$/caja 15 kilos'
$ws.Range("S32").Value = 1467
$ws.Range("T32").Value = 15

# Row 33
$ws.Range("D33").Value = 44179
$ws.Range("K33").Value = 'Modesto'
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 20000
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 20000
$ws.Range("S33").Value = 1111

# Row 34
$ws.Range("D34").Value = 44904
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 110
$ws.Range("N34").Value = 22000
$ws.Range("O34").Value = 22000
$ws.Range("P34").Value = 22000
$ws.Range("S34").Value = 1222

# Row 35
$ws.Range("D35").Value = 44546
$ws.Range("K35").Value = 'Castle Brite'
$ws.Range("M35").Value = 65
$ws.Range("N35").Value = 18000
$ws.Range("O35").Value = 18000
$ws.Range("P35").Value = 18000
$ws.Range("Q35").Value = '$/bandeja 18 kilos'
$ws.Range("R35").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S35").Value = 1000
$ws.Range("T35").Value = 18

# Row 36
$ws.Range("D36").Value = 44550
$ws.Range("K36").Value = 'Dina'
$ws.Range("M36").Value = 55
$ws.Range("N36").Value = 21000
$ws.Range("O36").Value = 21000
$ws.Range("P36").Value = 21000
$ws.Range("Q36").Value = '$/caja 18 kilos'
$ws.Range("R36").Value = 'Provincia de Quillota'
$ws.Range("S36").Value = 1167
$ws.Range("T36").Value = 18

# Row 37
$ws.Range("D37").Value = 44890
$ws.Range("K37").Value = 'Castle Brite'
$ws.Range("M37").Value = 80
$ws.Range("N37").Value = 20000
$ws.Range("O37").Value = 20000
$ws.Range("P37").Value = 20000
$ws.Range("R37").Value = 'Provincia de Limarí'
$ws.Range("S37").Value = 2000

# Row 38
$ws.Range("D38").Value = 44890
$ws.Range("M38").Value = 80
$ws.Range("O38").Value = 22000
$ws.Range("P38").Value = 22000
$ws.Range("Q38").Value = '$/caja 15 kilos'
$ws.Range("R38").Value = 'Provincia de Limarí'
$ws.Range("S38").Value = 1467
$ws.Range("T38").Value = 15

# Row 39
$ws.Range("D39").Value = 44176
$ws.Range("K39").Value = 'Modesto'
$ws.Range("M39").Value = 40
$ws.Range("O39").Value = 20000
$ws.Range("P39").Value = 20000
$ws.Range("Q39").Value = '$/caja 18 kilos'
$ws.Range("S39").Value = 1111

# Row 40
$ws.Range("D40").Value = 44545
$ws.Range("M40").Value = 85
$ws.Range("O40").Value = 18000
$ws.Range("P40").Value = 17588
$ws.Range("Q40").Value = '$/bandeja 18 kilos'
$ws.Range("R40").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S40").Value = 977
$ws.Range("T40").Value = 18

# Row 41
$ws.Range("D41").Value = 44188
$ws.Range("K41").Value = 'Modesto'
$ws.Range("L41").Value = 'Primera'
$ws.Range("M41").Value = 35
$ws.Range("N41").Value = 20000
$ws.Range("O41").Value = 20000
$ws.Range("P41").Value = 20000
$ws.Range("Q41").Value = '$/bandeja 18 kilos'
$ws.Range("R41").Value = 'Región de O''Higgins'
$ws.Range("S41").Value = 1111
$ws.Range("T41").Value = 18

# Row 42
$ws.Range("D42").Value = 44889
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 15000
$ws.Range("O42").Value = 15000
$ws.Range("P42").Value = 15000
$ws.Range("Q42").Value = '$/bandeja 10 kilos'
$ws.Range("R42").Value = 'Provincia de Limarí'
$ws.Range("S42").Value = 1500
$ws.Range("T42").Value = 10

# Row 43
$ws.Range("D43").Value = 44889
$ws.Range("M43").Value = 90
$ws.Range("N43").Value = 33000
$ws.Range("O43").Value = 34000
$ws.Range("P43").Value = 33556
$ws.Range("Q43").Value = '$/bandeja 18 kilos'
$ws.Range("R43").Value = 'Provincia de Limarí'
$ws.Range("S43").Value = 1864
$ws.Range("T43").Value = 18

# Row 44
$ws.Range("D44").Value = 44181
$ws.Range("K44").Value = 'Modesto'
$ws.Range("M44").Value = 140
$ws.Range("N44").Value = 17000
$ws.Range("O44").Value = 17000
$ws.Range("P44").Value = 17000
$ws.Range("Q44").Value = '$/caja 18 kilos'
$ws.Range("R44").Value = 'Región de O''Higgins'
$ws.Range("S44").Value = 944
$ws.Range("T44").Value = 18

# Row 45
$ws.Range("D45").Value = 44893
$ws.Range("K45").Value = 'Albaricoque'
$ws.Range("M45").Value = 250
$ws.Range("N45").Value = 15000
$ws.Range("O45").Value = 15000
$ws.Range("P45").Value = 15000
$ws.Range("Q45").Value = '$/bandeja 10 kilos'
$ws.Range("R45").Value = 'Provincia de Quillota'
$ws.Range("S45").Value = 1500
$ws.Range("T45").Value = 10

# Row 46
$ws.Range("D46").Value = 44175
$ws.Range("K46").Value = 'Castle Brite'
$ws.Range("M46").Value = 65
$ws.Range("Q46").Value = '$/bandeja 18 kilos'
$ws.Range("S46").Value = 1111
$ws.Range("T46").Value = 18

# Row 47
$ws.Range("D47").Value = 44175
$ws.Range("L47").Value = 'Segunda'
$ws.Range("M47").Value = 55
$ws.Range("N47").Value = 18000
$ws.Range("P47").Value = 18000
$ws.Range("R47").Value = 'Región de O''Higgins'
$ws.Range("S47").Value = 1000

# Row 48
$ws.Range("D48").Value = 44175
$ws.Range("L48").Value = 'Tercera'
$ws.Range("M48").Value = 45
$ws.Range("N48").Value = 14000
$ws.Range("O48").Value = 14000
$ws.Range("P48").Value = 14000
$ws.Range("Q48").Value = '$/bandeja 18 kilos'
$ws.Range("S48").Value = 778
$ws.Range("T48").Value = 18

# Row 49
$ws.Range("D49").Value = 44526
$ws.Range("M49").Value = 65
$ws.Range("Q49").Value = '$/bandeja 10 kilos'
$ws.Range("R49").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S49").Value = 2000
$ws.Range("T49").Value = 10

# Row 50
$ws.Range("D50").Value = 44202
$ws.Range("M50").Value = 25
$ws.Range("N50").Value = 18000
$ws.Range("O50").Value = 18000
$ws.Range("P50").Value = 18000
$ws.Range("Q50").Value = '$/caja 15 kilos'
$ws.Range("S50").Value = 1200
$ws.Range("T50").Value = 15

# Row 51
$ws.Range("D51").Value = 44160
$ws.Range("L51").Value = 'Primera'
$ws.Range("M51").Value = 175
$ws.Range("N51").Value = 18000
$ws.Range("O51").Value = 20000
$ws.Range("P51").Value = 18743
$ws.Range("Q51").Value = '$/caja 15 kilos'
$ws.Range("R51").Value = 'Provincia de Limarí'
$ws.Range("S51").Value = 1250
$ws.Range("T51").Value = 15

# Row 52
$ws.Range("D52").Value = 44166
$ws.Range("M52").Value = 400
$ws.Range("N52").Value = 17000
$ws.Range("O52").Value = 20000
$ws.Range("P52").Value = 17750
$ws.Range("Q52").Value = '$/caja 15 kilos'
$ws.Range("R52").Value = 'Provincia de Limarí'
$ws.Range("S52").Value = 1183
$ws.Range("T52").Value = 15

# Row 53
$ws.Range("D53").Value = 44166
$ws.Range("K53").Value = 'Castle Brite'
$ws.Range("L53").Value = 'Segunda'
$ws.Range("M53").Value = 80
$ws.Range("N53").Value = 12000
$ws.Range("O53").Value = 12000
$ws.Range("P53").Value = 12000
$ws.Range("Q53").Value = '$/caja 15 kilos'
$ws.Range("R53").Value = 'Provincia de Limarí'
$ws.Range("S53").Value = 800
$ws.Range("T53").Value = 15

# Row 54
$ws.Range("D54").Value = 44186
$ws.Range("K54").Value = 'Modesto'
$ws.Range("L54").Value = 'Primera'
$ws.Range("M54").Value = 55
$ws.Range("N54").Value = 20000
$ws.Range("O54").Value = 20000
$ws.Range("P54").Value = 20000
$ws.Range("Q54").Value = '$/bandeja 18 kilos'
$ws.Range("S54").Value = 1111
$ws.Range("T54").Value = 18

# Row 55
$ws.Range("D55").Value = 44917
$ws.Range("K55").Value = 'Modesto'
$ws.Range("M55").Value = 400
$ws.Range("N55").Value = 20000
$ws.Range("P55").Value = 20500
$ws.Range("S55").Value = 1139

# Row 56
$ws.Range("D56").Value = 44566
$ws.Range("K56").Value = 'Modesto'
$ws.Range("M56").Value = 55
$ws.Range("N56").Value = 15000
$ws.Range("O56").Value = 15000
$ws.Range("P56").Value = 15000
$ws.Range("Q56").Value = '$/bandeja 10 kilos'
$ws.Range("R56").Value = 'Provincia de Quillota'
$ws.Range("S56").Value = 1500
$ws.Range("T56").Value = 10

# Row 57
$ws.Range("D57").Value = 44167
$ws.Range("M57").Value = 100
$ws.Range("N57").Value = 20000
$ws.Range("O57").Value = 20000
$ws.Range("P57").Value = 20000
$ws.Range("Q57").Value = '$/caja 15 kilos'
$ws.Range("R57").Value = 'Región de O''Higgins'
$ws.Range("S57").Value = 1333
$ws.Range("T57").Value = 15

# Row 58
$ws.Range("D58").Value = 44217
$ws.Range("K58").Value = 'Modesto'
$ws.Range("L58").Value = 'Primera'
$ws.Range("N58").Value = 18000
$ws.Range("O58").Value = 18000
$ws.Range("P58").Value = 18000
$ws.Range("Q58").Value = '$/bandeja 18 kilos'
$ws.Range("R58").Value = 'Región de O''Higgins'
$ws.Range("S58").Value = 1000
$ws.Range("T58").Value = 18

# Row 59
$ws.Range("D59").Value = 44553
$ws.Range("K59").Value = 'Dina'
$ws.Range("M59").Value = 65
$ws.Range("N59").Value = 20000
$ws.Range("O59").Value = 20000
$ws.Range("P59").Value = 20000
$ws.Range("Q59").Value = '$/caja 18 kilos'
$ws.Range("R59").Value = 'Provincia de Quillota'
$ws.Range("S59").Value = 1111
$ws.Range("T59").Value = 18

# Row 60
$ws.Range("D60").Value = 44908
$ws.Range("M60").Value = 210
$ws.Range("N60").Value = 13000
$ws.Range("O60").Value = 13000
$ws.Range("P60").Value = 13000
$ws.Range("Q60").Value = '$/bandeja 10 kilos'
$ws.Range("R60").Value = 'Región de O''Higgins'
$ws.Range("S60").Value = 1300
$ws.Range("T60").Value = 10

# Row 61
$ws.Range("D61").Value = 44914
$ws.Range("L61").Value = 'Especial'
$ws.Range("M61").Value = 20
$ws.Range("N61").Value = 19000
$ws.Range("O61").Value = 19000
$ws.Range("P61").Value = 19000
$ws.Range("Q61").Value = '$/bandeja 10 kilos'
$ws.Range("R61").Value = 'Región de O''Higgins'
$ws.Range("S61").Value = 1900
$ws.Range("T61").Value = 10

# Row 62
$ws.Range("D62").Value = 44914
$ws.Range("L62").Value = 'Primera'
$ws.Range("M62").Value = 45
$ws.Range("N62").Value = 16000
$ws.Range("O62").Value = 16000
$ws.Range("P62").Value = 16000
$ws.Range("Q62").Value = '$/bandeja 10 kilos'
$ws.Range("R62").Value = 'Región de O''Higgins'
$ws.Range("S62").Value = 1600
$ws.Range("T62").Value = 10

# Row 63
$ws.Range("D63").Value = 44918
$ws.Range("K63").Value = 'Modesto'
$ws.Range("M63").Value = 200
$ws.Range("N63").Value = 20000
$ws.Range("O63").Value = 20000
$ws.Range("P63").Value = 20000
$ws.Range("Q63").Value = '$/bandeja 18 kilos'
$ws.Range("R63").Value = 'Región de O''Higgins'
$ws.Range("S63").Value = 1111

# Row 64
$ws.Range("D64").Value = 44529
$ws.Range("L64").Value = 'Primera'
$ws.Range("M64").Value = 75
$ws.Range("N64").Value = 20000
$ws.Range("O64").Value = 20000
$ws.Range("P64").Value = 20000
$ws.Range("Q64").Value = '$/bandeja 10 kilos'
$ws.Range("R64").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S64").Value = 2000
$ws.Range("T64").Value = 10

# Row 65
$ws.Range("D65").Value = 44529
$ws.Range("L65").Value = 'Segunda'
$ws.Range("M65").Value = 45
$ws.Range("N65").Value = 15000
$ws.Range("O65").Value = 15000
$ws.Range("P65").Value = 15000
$ws.Range("Q65").Value = '$/bandeja 10 kilos'
$ws.Range("R65").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S65").Value = 1500
$ws.Range("T65").Value = 10

# Row 66
$ws.Range("D66").Value = 44560
$ws.Range("M66").Value = 60
$ws.Range("N66").Value = 18000
$ws.Range("P66").Value = 18833
$ws.Range("Q66").Value = '$/caja 18 kilos'
$ws.Range("R66").Value = 'Provincia de Quillota'
$ws.Range("S66").Value = 1046
$ws.Range("T66").Value = 18

# Row 67
$ws.Range("D67").Value = 44524
$ws.Range("K67").Value = 'Castle Brite'
$ws.Range("M67").Value = 145
$ws.Range("N67").Value = 18000
$ws.Range("P67").Value = 19103
$ws.Range("Q67").Value = '$/bandeja 7 kilos'
$ws.Range("R67").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S67").Value = 2729
$ws.Range("T67").Value = 7

# Row 68
$ws.Range("D68").Value = 44901
$ws.Range("K68").Value = 'Castle Brite'
$ws.Range("M68").Value = 100
$ws.Range("N68").Value = 24000
$ws.Range("O68").Value = 24000
$ws.Range("P68").Value = 24000
$ws.Range("Q68").Value = '$/bandeja 18 kilos'
$ws.Range("R68").Value = 'Provincia de Limarí'
$ws.Range("S68").Value = 1333
$ws.Range("T68").Value = 18

# Row 69
$ws.Range("D69").Value = 44537
$ws.Range("K69").Value = 'Castle Brite'
$ws.Range("M69").Value = 115
$ws.Range("N69").Value = 18000
$ws.Range("P69").Value = 19130
$ws.Range("Q69").Value = '$/bandeja 18 kilos'
$ws.Range("R69").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S69").Value = 1063

# Row 70
$ws.Range("D70").Value = 44922
$ws.Range("K70").Value = 'Dina'
$ws.Range("M70").Value = 35
$ws.Range("N70").Value = 20000
$ws.Range("O70").Value = 20000
$ws.Range("P70").Value = 20000
$ws.Range("S70").Value = 1333

# Row 71
$ws.Range("D71").Value = 44897
$ws.Range("L71").Value = 'Especial'
$ws.Range("M71").Value = 95
$ws.Range("N71").Value = 18000
$ws.Range("O71").Value = 18000
$ws.Range("P71").Value = 18000
$ws.Range("Q71").Value = '$/bandeja 7 kilos'
$ws.Range("S71").Value = 2571
$ws.Range("T71").Value = 7

# Row 72
$ws.Range("D72").Value = 44897
$ws.Range("K72").Value = 'Castle Brite'
$ws.Range("M72").Value = 45
$ws.Range("N72").Value = 28000
$ws.Range("O72").Value = 28000
$ws.Range("P72").Value = 28000
$ws.Range("Q72").Value = '$/bandeja 18 kilos'
$ws.Range("R72").Value = 'Provincia de Limarí'
$ws.Range("S72").Value = 1556
$ws.Range("T72").Value = 18

# Row 73
$ws.Range("D73").Value = 44544
$ws.Range("M73").Value = 35
$ws.Range("O73").Value = 18000
$ws.Range("P73").Value = 18000
$ws.Range("S73").Value = 1000

# Row 74
$ws.Range("D74").Value = 44895
$ws.Range("N74").Value = 28000
$ws.Range("O74").Value = 28000
$ws.Range("P74").Value = 28000
$ws.Range("Q74").Value = '$/bandeja 18 kilos'
$ws.Range("R74").Value = 'Provincia de Limarí'
$ws.Range("S74").Value = 1556
$ws.Range("T74").Value = 18

# Row 75
$ws.Range("D75").Value = 44522
$ws.Range("M75").Value = 30
$ws.Range("N75").Value = 35000
$ws.Range("O75").Value = 35000
$ws.Range("P75").Value = 35000
$ws.Range("Q75").Value = '$/caja 15 kilos'
$ws.Range("R75").Value = 'Provincia de Limarí'
$ws.Range("S75").Value = 2333
$ws.Range("T75").Value = 15

# Row 76
$ws.Range("D76").Value = 44900
$ws.Range("L76").Value = 'Especial'
$ws.Range("M76").Value = 50
$ws.Range("N76").Value = 27000
$ws.Range("O76").Value = 27000
$ws.Range("P76").Value = 27000
$ws.Range("S76").Value = 1500

# Row 77
$ws.Range("D77").Value = 44900
$ws.Range("L77").Value = 'Especial'
$ws.Range("M77").Value = 200
$ws.Range("N77").Value = 22000
$ws.Range("O77").Value = 22000
$ws.Range("P77").Value = 22000
$ws.Range("Q77").Value = '$/caja 18 kilos'
$ws.Range("S77").Value = 1222

# Row 78
$ws.Range("D78").Value = 44900
$ws.Range("K78").Value = 'Castle Brite'
$ws.Range("M78").Value = 480
$ws.Range("N78").Value = 20000
$ws.Range("O78").Value = 25000
$ws.Range("P78").Value = 23062
$ws.Range("Q78").Value = '$/bandeja 18 kilos'
$ws.Range("S78").Value = 1281

# Row 79
$ws.Range("D79").Value = 44900
$ws.Range("M79").Value = 250
$ws.Range("N79").Value = 22000
$ws.Range("O79").Value = 22000
$ws.Range("P79").Value = 22000
$ws.Range("Q79").Value = '$/caja 18 kilos'
$ws.Range("R79").Value = 'Región de O''Higgins'
$ws.Range("S79").Value = 1222

# Row 80
$ws.Range("D80").Value = 44900
$ws.Range("L80").Value = 'Segunda'
$ws.Range("N80").Value = 20000
$ws.Range("O80").Value = 20000
$ws.Range("P80").Value = 20000
$ws.Range("R80").Value = 'Región de O''Higgins'
$ws.Range("S80").Value = 1111

# Row 81
$ws.Range("D81").Value = 44532
$ws.Range("K81").Value = 'Castle Brite'
$ws.Range("L81").Value = 'Especial'
$ws.Range("N81").Value = 20000
$ws.Range("O81").Value = 20000
$ws.Range("P81").Value = 20000
$ws.Range("Q81").Value = '$/bandeja 10 kilos'
$ws.Range("R81").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S81").Value = 2000
$ws.Range("T81").Value = 10

# Row 82
$ws.Range("D82").Value = 44532
$ws.Range("K82").Value = 'Castle Brite'
$ws.Range("M82").Value = 420
$ws.Range("N82").Value = 13000
$ws.Range("P82").Value = 14048
$ws.Range("R82").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S82").Value = 1405

# Row 83
$ws.Range("D83").Value = 44532
$ws.Range("K83").Value = 'Castle Brite'
$ws.Range("M83").Value = 2
$ws.Range("N83").Value = 600000
$ws.Range("O83").Value = 600000
$ws.Range("P83").Value = 600000
$ws.Range("Q83").Value = '$/bins (500 kilos)'
$ws.Range("S83").Value = 1200
$ws.Range("T83").Value = 500

# Row 84
$ws.Range("D84").Value = 44532
$ws.Range("K84").Value = 'Castle Brite'
$ws.Range("M84").Value = 100
$ws.Range("N84").Value = 30000
$ws.Range("O84").Value = 30000
$ws.Range("P84").Value = 30000
$ws.Range("R84").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S84").Value = 1667

# Row 85
$ws.Range("D85").Value = 44532
$ws.Range("L85").Value = 'Segunda'
$ws.Range("M85").Value = 80
$ws.Range("N85").Value = 15000
$ws.Range("O85").Value = 15000
$ws.Range("P85").Value = 15000
$ws.Range("Q85").Value = '$/bandeja 18 kilos'
$ws.Range("S85").Value = 833

# Row 86
$ws.Range("D86").Value = 44201
$ws.Range("K86").Value = 'Modesto'
$ws.Range("M86").Value = 45
$ws.Range("N86").Value = 18000
$ws.Range("O86").Value = 18000
$ws.Range("P86").Value = 18000
$ws.Range("Q86").Value = '$/caja 15 kilos'
$ws.Range("R86").Value = 'Región de O''Higgins'
$ws.Range("S86").Value = 1200
$ws.Range("T86").Value = 15

# Row 87
$ws.Range("D87").Value = 44567
$ws.Range("K87").Value = 'Modesto'
$ws.Range("L87").Value = 'Primera'
$ws.Range("M87").Value = 25
$ws.Range("N87").Value = 15000
$ws.Range("O87").Value = 15000
$ws.Range("P87").Value = 15000
$ws.Range("R87").Value = 'Provincia de Quillota'
$ws.Range("S87").Value = 1500

# Row 88
$ws.Range("D88").Value = 44525
$ws.Range("L88").Value = 'Primera'
$ws.Range("M88").Value = 55
$ws.Range("N88").Value = 20000
$ws.Range("O88").Value = 20000
$ws.Range("P88").Value = 20000
$ws.Range("S88").Value = 2000

# Row 89
$ws.Range("D89").Value = 44557
$ws.Range("K89").Value = 'Dina'
$ws.Range("M89").Value = 95
$ws.Range("N89").Value = 7000
$ws.Range("O89").Value = 7000
$ws.Range("P89").Value = 7000
$ws.Range("Q89").Value = '$/bandeja 6 kilos'
$ws.Range("R89").Value = 'Provincia de San Felipe de Aconcagua'
$ws.Range("S89").Value = 1167
$ws.Range("T89").Value = 6

# Row 90
$ws.Range("D90").Value = 44557
$ws.Range("K90").Value = 'Dina'
$ws.Range("M90").Value = 35
$ws.Range("N90").Value = 20000
$ws.Range("O90").Value = 20000
$ws.Range("P90").Value = 20000
$ws.Range("Q90").Value = '$/caja 18 kilos'
$ws.Range("R90").Value = 'Provincia de Quillota'
$ws.Range("S90").Value = 1111
$ws.Range("T90").Value = 18
